# Weekly CompStat update: new crime data collected (101st Precinct, week ending 3/24/2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump the report "Number" (Volume 31, Number 11 -> Number 12) ---
# A8 rich text: "Volume 31   Number  " + "11" -> "12"
$hdrNumber = $ws.Range("A8")
$hdrText = $hdrNumber.Value()
$idx = $hdrText.LastIndexOf("11") + 1
$hdrNumber.Characters($idx, 2).Text = "12"

# --- Header: advance the reporting week (3/11/2024-3/17/2024 -> 3/18/2024-3/24/2024) ---
$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 9).Text = "3/18/2024"
$weekCell.Characters(47, 9).Text = "3/24/2024"

# --- Crime-complaint table refresh (rows 14-30): new weekly/28-day/YTD/2-year figures ---
$ws.Range("N14").Value = -75
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 2
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 20
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = 3.571428571428
$ws.Range("L16").Value = 81.25
$ws.Range("M16").Value = -14.705882352941
$ws.Range("N16").Value = -77.165354330708
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 172.727272727273
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = 78.048780487804
$ws.Range("L17").Value = 37.735849056603
$ws.Range("M17").Value = 87.179487179487
$ws.Range("N17").Value = -21.505376344086
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -5.263157894736
$ws.Range("L18").Value = -10
$ws.Range("M18").Value = -28
$ws.Range("N18").Value = -84.873949579831
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -87.5
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = -50
$ws.Range("J19").Value = 42
$ws.Range("K19").Value = -30.952380952381
$ws.Range("L19").Value = -30.952380952381
$ws.Range("N19").Value = -50.847457627118
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = "0"
$ws.Range("H20").Value = "***.*"
$ws.Range("I20").Value = 15
$ws.Range("K20").Value = -6.25
$ws.Range("L20").Value = 114.285714285714
$ws.Range("M20").Value = -28.571428571428
$ws.Range("N20").Value = -86.725663716814
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -23.529411764705
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = 40.90909090909
$ws.Range("I21").Value = 171
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 14
$ws.Range("L21").Value = 18.75
$ws.Range("M21").Value = 17.931034482758
$ws.Range("N21").Value = -67.178502879078
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 20
$ws.Range("I23").Value = 19
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 11.764705882352
$ws.Range("M23").Value = 58.333333333333
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -70
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 45
$ws.Range("H24").Value = -8.888888888888
$ws.Range("I24").Value = 136
$ws.Range("J24").Value = 116
$ws.Range("K24").Value = 17.241379310344
$ws.Range("L24").Value = 23.636363636363
$ws.Range("M24").Value = 91.549295774647
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 4
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 10
$ws.Range("J25").Value = 15
$ws.Range("K25").Value = -33.333333333333
$ws.Range("L25").Value = 0
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -61.538461538461
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 81
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 1.25
$ws.Range("L26").Value = 3.846153846153
$ws.Range("M26").Value = 0
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -10
$ws.Range("L27").Value = 28.571428571428
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -66.666666666666
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = -20
$ws.Range("L29").Value = -85.714285714285
$ws.Range("M29").Value = -75
$ws.Range("N29").Value = -91.666666666666
$ws.Range("L30").Value = -80
$ws.Range("M30").Value = -75
$ws.Range("N30").Value = -91.666666666666
